$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, clone the formatting of row 46 (style + row height) into the new row 47,
# so the appended requisito line matches the rest of the "Requisitos" list.
$ws.Range("B46:C46").Copy()
$ws.Range("B47:C47").PasteSpecial(-4122)
$ws.Rows.Item(47).RowHeight = $ws.Rows.Item(46).RowHeight

# Requisito list changed (LOQ4031/LOQ4073 replaced by LOQ4010, LOQ4095 kept,
# two new requisitos LOQ4097/LOQ4098 added at the end):
#   row44: LOQ4031 -> LOQ4010
#   row45: LOQ4073 -> LOQ4095
#   row46: LOQ4095 -> LOQ4097
#   row47: (new)   -> LOQ4098
$ws.Range("B44").Value = "LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)`n"
$ws.Range("C44").Value = "LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)`n"

$ws.Range("B45").Value = "LOQ4095 -  Química Geral Experimental  (Requisito)`n"
$ws.Range("C45").Value = "LOQ4095 -  Química Geral Experimental  (Requisito)`n"

$ws.Range("B46").Value = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)`n"
$ws.Range("C46").Value = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)`n"

$ws.Range("B47").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)`n"
$ws.Range("C47").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)`n"
